$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 32258414  # H33: 33333694 -> 32258414
$ws.Cells.Item(33, 9).Value = 33333684  # I33: 34483120 -> 33333684
$ws.Cells.Item(33, 11).Value = 33333684  # K33: 34483120 -> 33333684
$ws.Cells.Item(33, 13).Value = -33333455  # M33: -34482891 -> -33333455

# Hunk 1: ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1615.3704  # H40: 1628.2593 -> 1615.3704
$ws.Cells.Item(40, 9).Value = 1466.375  # I40: 1459.5294 -> 1466.375
$ws.Cells.Item(40, 10).Value = 1832.091  # J40: 1915.1 -> 1832.091
$ws.Cells.Item(40, 11).Value = 1466.375  # K40: 1459.5294 -> 1466.375
$ws.Cells.Item(40, 12).Value = 1832.091  # L40: 1915.1 -> 1832.091
$ws.Cells.Item(40, 13).Value = -1291.375  # M40: -1284.5294 -> -1291.375
$ws.Cells.Item(40, 14).Value = -2182.091  # N40: -2265.1 -> -2182.091

# Hunk 2: ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2583.3333  # H62: 2371.2856 -> 2583.3333
$ws.Cells.Item(62, 9).Value = 5450  # I62: 5500 -> 5450
$ws.Cells.Item(62, 10).Value = 1150  # J62: 1119.8 -> 1150
$ws.Cells.Item(62, 11).Value = 5450  # K62: 5500 -> 5450
$ws.Cells.Item(62, 12).Value = 1150  # L62: 1119.8 -> 1150
$ws.Cells.Item(62, 13).Value = -4826  # M62: -4876 -> -4826
$ws.Cells.Item(62, 14).Value = -2398  # N62: -2367.8 -> -2398

# Hunk 3: ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2583.3333  # H65: 2371.2856 -> 2583.3333
$ws.Cells.Item(65, 9).Value = 5450  # I65: 5500 -> 5450
$ws.Cells.Item(65, 10).Value = 1150  # J65: 1119.8 -> 1150
$ws.Cells.Item(65, 11).Value = 27250  # K65: 27500 -> 27250
$ws.Cells.Item(65, 12).Value = 5750  # L65: 5599 -> 5750
$ws.Cells.Item(65, 13).Value = -24130  # M65: -24380 -> -24130
$ws.Cells.Item(65, 14).Value = -11990  # N65: -11839 -> -11990

# Hunk 4: ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 1818.5  # H92: 1785.76 -> 1818.5
$ws.Cells.Item(92, 9).Value = 1677.2  # I92: 1602.0952 -> 1677.2
$ws.Cells.Item(92, 10).Value = 2525  # J92: 2750 -> 2525
$ws.Cells.Item(92, 11).Value = 1677.2  # K92: 1602.0952 -> 1677.2
$ws.Cells.Item(92, 12).Value = 2525  # L92: 2750 -> 2525
$ws.Cells.Item(92, 13).Value = -429.2  # M92: -354.0952 -> -429.2
$ws.Cells.Item(92, 14).Value = -5021  # N92: -5246 -> -5021

# Hunk 5: ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 3829  # H94: 3786.0908 -> 3829
$ws.Cells.Item(94, 9).Value = 1985.7142  # I94: 1129.4 -> 1985.7142
$ws.Cells.Item(94, 10).Value = 5672.2856  # J94: 6000 -> 5672.2856
$ws.Cells.Item(94, 11).Value = 1985.7142  # K94: 1129.4 -> 1985.7142
$ws.Cells.Item(94, 12).Value = 5672.2856  # L94: 6000 -> 5672.2856
$ws.Cells.Item(94, 13).Value = -1534.7142  # M94: -678.4000000000001 -> -1534.7142
$ws.Cells.Item(94, 14).Value = -6574.2856  # N94: -6902 -> -6574.2856

# Hunk 6: ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 1336.3529  # H96: 1338.1177 -> 1336.3529
$ws.Cells.Item(96, 9).Value = 2756.5  # I96: 3575.3333 -> 2756.5
$ws.Cells.Item(96, 10).Value = 899.38464  # J96: 858.7143 -> 899.38464
$ws.Cells.Item(96, 11).Value = 8269.5  # K96: 10725.9999 -> 8269.5
$ws.Cells.Item(96, 12).Value = 2698.15392  # L96: 2576.1429 -> 2698.15392
$ws.Cells.Item(96, 13).Value = -6896.5  # M96: -9352.999899999999 -> -6896.5
$ws.Cells.Item(96, 14).Value = -5444.15392  # N96: -5322.1429 -> -5444.15392

# Hunk 7: ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(99, 8).Value = 78374.46000000001  # H99: 72897.71000000001 -> 78374.46000000001
$ws.Cells.Item(99, 9).Value = 662.25  # I99: 699.75 -> 662.25
$ws.Cells.Item(99, 10).Value = 202714  # J99: 169161.67 -> 202714
$ws.Cells.Item(99, 11).Value = 1986.75  # K99: 2099.25 -> 1986.75
$ws.Cells.Item(99, 12).Value = 608142  # L99: 507485.01 -> 608142
$ws.Cells.Item(99, 13).Value = -488.75  # M99: -601.25 -> -488.75
$ws.Cells.Item(99, 14).Value = -611138  # N99: -510481.01 -> -611138

# Hunk 8: ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 3069.077  # H100: 3263.6365 -> 3069.077
$ws.Cells.Item(100, 9).Value = 2499.75  # I100: 2666.6667 -> 2499.75
$ws.Cells.Item(100, 10).Value = 3322.111  # J100: 3487.5 -> 3322.111
$ws.Cells.Item(100, 11).Value = 2499.75  # K100: 2666.6667 -> 2499.75
$ws.Cells.Item(100, 12).Value = 3322.111  # L100: 3487.5 -> 3322.111
$ws.Cells.Item(100, 13).Value = -1958.75  # M100: -2125.6667 -> -1958.75
$ws.Cells.Item(100, 14).Value = -4404.111  # N100: -4569.5 -> -4404.111

# Hunk 9: ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 2229.475  # H129: 2116.7954 -> 2229.475
$ws.Cells.Item(129, 10).Value = 2441.4285  # J129: 2292.5642 -> 2441.4285
$ws.Cells.Item(129, 12).Value = 7324.2855  # L129: 6877.692599999999 -> 7324.2855
$ws.Cells.Item(129, 14).Value = -17324.2855  # N129: -16877.6926 -> -17324.2855

# Hunk 10: ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 2771.2563  # H132: 3499.9688 -> 2771.2563
$ws.Cells.Item(132, 9).Value = 2554.7104  # I132: 3258.0322 -> 2554.7104
$ws.Cells.Item(132, 11).Value = 7664.1312  # K132: 9774.096600000001 -> 7664.1312
$ws.Cells.Item(132, 13).Value = -5134.1312  # M132: -7244.096600000001 -> -5134.1312

# Hunk 11: ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 17859492  # H135: 17859508 -> 17859492
$ws.Cells.Item(135, 9).Value = 1892.8422  # I135: 1894.421 -> 1892.8422
$ws.Cells.Item(135, 10).Value = 55558868  # J135: 55558910 -> 55558868
$ws.Cells.Item(135, 11).Value = 17035.5798  # K135: 17049.789 -> 17035.5798
$ws.Cells.Item(135, 12).Value = 500029812  # L135: 500030190 -> 500029812
$ws.Cells.Item(135, 13).Value = -14500.5798  # M135: -14514.789 -> -14500.5798
$ws.Cells.Item(135, 14).Value = -500034882  # N135: -500035260 -> -500034882

# Hunk 12: ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(136, 8).Value = 49975  # H136: 50000 -> 49975
$ws.Cells.Item(136, 10).Value = 49975  # J136: 50000 -> 49975
$ws.Cells.Item(136, 12).Value = 49975  # L136: 50000 -> 49975
$ws.Cells.Item(136, 14).Value = -60175  # N136: -60200 -> -60175

# Hunk 13: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 16668768  # H137: 14287759 -> 16668768
$ws.Cells.Item(137, 9).Value = 1150.3334  # I137: 1466.6666 -> 1150.3334
$ws.Cells.Item(137, 10).Value = 33336384  # J137: 18184020 -> 33336384
$ws.Cells.Item(137, 11).Value = 3451.0002  # K137: 4399.9998 -> 3451.0002
$ws.Cells.Item(137, 12).Value = 100009152  # L137: 54552060 -> 100009152
$ws.Cells.Item(137, 13).Value = -901.0001999999999  # M137: -1849.9998 -> -901.0001999999999
$ws.Cells.Item(137, 14).Value = -100014252  # N137: -54557160 -> -100014252

# Hunk 14: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2625.9033  # H138: 2776.8667 -> 2625.9033
$ws.Cells.Item(138, 9).Value = 1900.3334  # I138: 1993.4286 -> 1900.3334
$ws.Cells.Item(138, 10).Value = 3306.125  # J138: 3462.375 -> 3306.125
$ws.Cells.Item(138, 11).Value = 5701.0002  # K138: 5980.2858 -> 5701.0002
$ws.Cells.Item(138, 12).Value = 9918.375  # L138: 10387.125 -> 9918.375
$ws.Cells.Item(138, 13).Value = -561.0002000000004  # M138: -840.2857999999997 -> -561.0002000000004
$ws.Cells.Item(138, 14).Value = -20198.375  # N138: -20667.125 -> -20198.375

# Hunk 15: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 125006260  # H61: 12501798 -> 125006260
$ws.Cells.Item(61, 9).Value = 125006260  # I61: 13515268 -> 125006260
$ws.Cells.Item(61, 10).Value = 0  # J61: 2342.6667 -> 0
$ws.Cells.Item(61, 11).Value = 125006260  # K61: 13515268 -> 125006260
$ws.Cells.Item(61, 12).Value = 0  # L61: 2342.6667 -> 0
$ws.Cells.Item(61, 13).Value = -125006048  # M61: -13515056 -> -125006048
$ws.Cells.Item(61, 14).ClearContents()  # N61: -2766.6667 -> (removed)

# Hunk 16: ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 6077.0557  # H97: 6492.2354 -> 6077.0557
$ws.Cells.Item(97, 9).Value = 7612.5713  # I97: 8273.615 -> 7612.5713
$ws.Cells.Item(97, 11).Value = 7612.5713  # K97: 8273.615 -> 7612.5713
$ws.Cells.Item(97, 13).Value = -7116.5713  # M97: -7777.615 -> -7116.5713

# Hunk 17: ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1361.4814  # H102: 1566.25 -> 1361.4814
$ws.Cells.Item(102, 9).Value = 1030.4  # I102: 1208.6364 -> 1030.4
$ws.Cells.Item(102, 11).Value = 1030.4  # K102: 1208.6364 -> 1030.4
$ws.Cells.Item(102, 13).Value = 591.5999999999999  # M102: 413.3635999999999 -> 591.5999999999999

# Hunk 18: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 125006260  # H136: 12501798 -> 125006260
$ws.Cells.Item(136, 9).Value = 125006260  # I136: 13515268 -> 125006260
$ws.Cells.Item(136, 10).Value = 0  # J136: 2342.6667 -> 0
$ws.Cells.Item(136, 11).Value = 375018780  # K136: 40545804 -> 375018780
$ws.Cells.Item(136, 12).Value = 0  # L136: 7028.000100000001 -> 0
$ws.Cells.Item(136, 13).Value = -375016230  # M136: -40543254 -> -375016230
$ws.Cells.Item(136, 14).ClearContents()  # N136: -12128.0001 -> (removed)

# Hunk 19: ARM row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(137, 8).Value = 64976.668  # H137: 74997.78 -> 64976.668
$ws.Cells.Item(137, 10).Value = 79965  # J137: 79997.5 -> 79965
$ws.Cells.Item(137, 12).Value = 79965  # L137: 79997.5 -> 79965
$ws.Cells.Item(137, 14).Value = -90165  # N137: -90197.5 -> -90165

# Hunk 20: BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 33335312  # H86: 33335306 -> 33335312
$ws.Cells.Item(86, 9).Value = 1968.5714  # I86: 1960 -> 1968.5714
$ws.Cells.Item(86, 10).Value = 62501988  # J86: 71430560 -> 62501988
$ws.Cells.Item(86, 11).Value = 1968.5714  # K86: 1960 -> 1968.5714
$ws.Cells.Item(86, 12).Value = 62501988  # L86: 71430560 -> 62501988
$ws.Cells.Item(86, 13).Value = -845.5714  # M86: -837 -> -845.5714
$ws.Cells.Item(86, 14).Value = -62504234  # N86: -71432806 -> -62504234

# Hunk 21: BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 33335312  # H89: 33335306 -> 33335312
$ws.Cells.Item(89, 9).Value = 1968.5714  # I89: 1960 -> 1968.5714
$ws.Cells.Item(89, 10).Value = 62501988  # J89: 71430560 -> 62501988
$ws.Cells.Item(89, 11).Value = 9842.857  # K89: 9800 -> 9842.857
$ws.Cells.Item(89, 12).Value = 312509940  # L89: 357152800 -> 312509940
$ws.Cells.Item(89, 13).Value = -4226.857  # M89: -4184 -> -4226.857
$ws.Cells.Item(89, 14).Value = -312521172  # N89: -357164032 -> -312521172

# Hunk 22: CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 5473.2  # H58: 4281 -> 5473.2
$ws.Cells.Item(58, 9).Value = 788.6667  # I58: 749.6667 -> 788.6667
$ws.Cells.Item(58, 10).Value = 12500  # J58: 6399.8 -> 12500
$ws.Cells.Item(58, 11).Value = 788.6667  # K58: 749.6667 -> 788.6667
$ws.Cells.Item(58, 12).Value = 12500  # L58: 6399.8 -> 12500
$ws.Cells.Item(58, 13).Value = -585.6667  # M58: -546.6667 -> -585.6667
$ws.Cells.Item(58, 14).Value = -12906  # N58: -6805.8 -> -12906

# Hunk 23: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3610.5  # H132: 3190 -> 3610.5
$ws.Cells.Item(132, 9).Value = 2682.4375  # I132: 2321.4 -> 2682.4375
$ws.Cells.Item(132, 11).Value = 8047.3125  # K132: 6964.200000000001 -> 8047.3125
$ws.Cells.Item(132, 13).Value = -5517.3125  # M132: -4434.200000000001 -> -5517.3125

# Hunk 24: CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 5473.2  # H136: 4281 -> 5473.2
$ws.Cells.Item(136, 9).Value = 788.6667  # I136: 749.6667 -> 788.6667
$ws.Cells.Item(136, 10).Value = 12500  # J136: 6399.8 -> 12500
$ws.Cells.Item(136, 11).Value = 2366.0001  # K136: 2249.0001 -> 2366.0001
$ws.Cells.Item(136, 12).Value = 37500  # L136: 19199.4 -> 37500
$ws.Cells.Item(136, 13).Value = 183.9998999999998  # M136: 300.9998999999998 -> 183.9998999999998
$ws.Cells.Item(136, 14).Value = -42600  # N136: -24299.4 -> -42600

# Hunk 25: CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(140, 8).Value = 33068.57  # H140: 31158 -> 33068.57
$ws.Cells.Item(140, 10).Value = 33068.57  # J140: 31158 -> 33068.57
$ws.Cells.Item(140, 12).Value = 33068.57  # L140: 31158 -> 33068.57
$ws.Cells.Item(140, 14).Value = -43428.57  # N140: -41518 -> -43428.57

# Hunk 26: CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 69920  # H37: 75233.336 -> 69920
$ws.Cells.Item(37, 10).Value = 69920  # J37: 75233.336 -> 69920
$ws.Cells.Item(37, 12).Value = 209760  # L37: 225700.008 -> 209760
$ws.Cells.Item(37, 14).Value = -209984  # N37: -225924.008 -> -209984

# Hunk 27: CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 14953.846  # H87: 12253.866 -> 14953.846
$ws.Cells.Item(87, 9).Value = 10225  # I87: 5051.3335 -> 10225
$ws.Cells.Item(87, 11).Value = 30675  # K87: 15154.0005 -> 30675
$ws.Cells.Item(87, 13).Value = -29427  # M87: -13906.0005 -> -29427

# Hunk 28: CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 14953.846  # H90: 12253.866 -> 14953.846
$ws.Cells.Item(90, 9).Value = 10225  # I90: 5051.3335 -> 10225
$ws.Cells.Item(90, 11).Value = 92025  # K90: 45462.0015 -> 92025
$ws.Cells.Item(90, 13).Value = -85785  # M90: -39222.0015 -> -85785

# Hunk 29: CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1360.9546  # H113: 1356.8636 -> 1360.9546
$ws.Cells.Item(113, 10).Value = 2097.6667  # J113: 2090.1667 -> 2097.6667
$ws.Cells.Item(113, 12).Value = 6293.000100000001  # L113: 6270.500100000001 -> 6293.000100000001
$ws.Cells.Item(113, 14).Value = -10633.0001  # N113: -10610.5001 -> -10633.0001

# Hunk 30: CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1005.0571  # H122: 1073.2188 -> 1005.0571
$ws.Cells.Item(122, 9).Value = 829.25  # I122: 908 -> 829.25
$ws.Cells.Item(122, 11).Value = 7463.25  # K122: 8172 -> 7463.25
$ws.Cells.Item(122, 13).Value = -5013.25  # M122: -5722 -> -5013.25

# Hunk 31: GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 44344.6  # H70: 70876.55499999999 -> 44344.6
$ws.Cells.Item(70, 9).Value = 123246  # I70: 202660 -> 123246
$ws.Cells.Item(70, 10).Value = 4893.9  # J70: 4984.8335 -> 4893.9
$ws.Cells.Item(70, 11).Value = 123246  # K70: 202660 -> 123246
$ws.Cells.Item(70, 12).Value = 4893.9  # L70: 4984.8335 -> 4893.9
$ws.Cells.Item(70, 13).Value = -122976  # M70: -202390 -> -122976
$ws.Cells.Item(70, 14).Value = -5433.9  # N70: -5524.8335 -> -5433.9

# Hunk 32: GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 44344.6  # H73: 70876.55499999999 -> 44344.6
$ws.Cells.Item(73, 9).Value = 123246  # I73: 202660 -> 123246
$ws.Cells.Item(73, 10).Value = 4893.9  # J73: 4984.8335 -> 4893.9
$ws.Cells.Item(73, 11).Value = 123246  # K73: 202660 -> 123246
$ws.Cells.Item(73, 12).Value = 4893.9  # L73: 4984.8335 -> 4893.9
$ws.Cells.Item(73, 13).Value = -122310  # M73: -201724 -> -122310
$ws.Cells.Item(73, 14).Value = -6765.9  # N73: -6856.8335 -> -6765.9

# Hunk 33: GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 922.3871  # H97: 929.4838999999999 -> 922.3871
$ws.Cells.Item(97, 9).Value = 1041.381  # I97: 1045.1904 -> 1041.381
$ws.Cells.Item(97, 10).Value = 672.5  # J97: 686.5 -> 672.5
$ws.Cells.Item(97, 11).Value = 1041.381  # K97: 1045.1904 -> 1041.381
$ws.Cells.Item(97, 12).Value = 672.5  # L97: 686.5 -> 672.5
$ws.Cells.Item(97, 13).Value = -545.3810000000001  # M97: -549.1904 -> -545.3810000000001
$ws.Cells.Item(97, 14).Value = -1664.5  # N97: -1678.5 -> -1664.5

# Hunk 34: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4007.139  # H132: 3465.1333 -> 4007.139
$ws.Cells.Item(132, 9).Value = 2862.625  # I132: 2340.3438 -> 2862.625
$ws.Cells.Item(132, 10).Value = 6296.1665  # J132: 6233.846 -> 6296.1665
$ws.Cells.Item(132, 11).Value = 8587.875  # K132: 7021.0314 -> 8587.875
$ws.Cells.Item(132, 12).Value = 18888.4995  # L132: 18701.538 -> 18888.4995
$ws.Cells.Item(132, 13).Value = -6057.875  # M132: -4491.0314 -> -6057.875
$ws.Cells.Item(132, 14).Value = -23948.4995  # N132: -23761.538 -> -23948.4995

# Hunk 35: GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 57665.332  # H138: 57656 -> 57665.332
$ws.Cells.Item(138, 10).Value = 57665.332  # J138: 57656 -> 57665.332
$ws.Cells.Item(138, 12).Value = 57665.332  # L138: 57656 -> 57665.332
$ws.Cells.Item(138, 14).Value = -67945.33199999999  # N138: -67936 -> -67945.33199999999

# Hunk 36: LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4605.2705  # H7: 4599.7295 -> 4605.2705
$ws.Cells.Item(7, 9).Value = 4505.5884  # I7: 4493.5293 -> 4505.5884
$ws.Cells.Item(7, 11).Value = 4505.5884  # K7: 4493.5293 -> 4505.5884
$ws.Cells.Item(7, 13).Value = -4393.5884  # M7: -4381.5293 -> -4393.5884

# Hunk 37: LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3793.394  # H40: 5363.6 -> 3793.394
$ws.Cells.Item(40, 9).Value = 5291.2856  # I40: 5291 -> 5291.2856
$ws.Cells.Item(40, 10).Value = 2689.6843  # J40: 5533 -> 2689.6843
$ws.Cells.Item(40, 11).Value = 5291.2856  # K40: 5291 -> 5291.2856
$ws.Cells.Item(40, 12).Value = 2689.6843  # L40: 5533 -> 2689.6843
$ws.Cells.Item(40, 13).Value = -5155.2856  # M40: -5155 -> -5155.2856
$ws.Cells.Item(40, 14).Value = -2961.6843  # N40: -5805 -> -2961.6843

# Hunk 38: LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2350  # H93: 1468.75 -> 2350
$ws.Cells.Item(93, 9).Value = 2525  # I93: 1389.2222 -> 2525
$ws.Cells.Item(93, 10).Value = 2000  # J93: 1571 -> 2000
$ws.Cells.Item(93, 11).Value = 2525  # K93: 1389.2222 -> 2525
$ws.Cells.Item(93, 12).Value = 2000  # L93: 1571 -> 2000
$ws.Cells.Item(93, 13).Value = -1277  # M93: -141.2221999999999 -> -1277
$ws.Cells.Item(93, 14).Value = -4496  # N93: -4067 -> -4496

# Hunk 39: LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 4605.2705  # H126: 4599.7295 -> 4605.2705
$ws.Cells.Item(126, 9).Value = 4505.5884  # I126: 4493.5293 -> 4505.5884
$ws.Cells.Item(126, 11).Value = 13516.7652  # K126: 13480.5879 -> 13516.7652
$ws.Cells.Item(126, 13).Value = -11046.7652  # M126: -11010.5879 -> -11046.7652

# Hunk 40: LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 27787532  # H136: 23817978 -> 27787532
$ws.Cells.Item(136, 9).Value = 35715612  # I136: 29412972 -> 35715612
$ws.Cells.Item(136, 11).Value = 107146836  # K136: 88238916 -> 107146836
$ws.Cells.Item(136, 13).Value = -107144286  # M136: -88236366 -> -107144286

# Hunk 41: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1053.0161  # H132: 1071 -> 1053.0161
$ws.Cells.Item(132, 9).Value = 740.73334  # I132: 749.4888999999999 -> 740.73334
$ws.Cells.Item(132, 10).Value = 1879.6471  # J132: 1975.25 -> 1879.6471
$ws.Cells.Item(132, 11).Value = 2222.20002  # K132: 2248.4667 -> 2222.20002
$ws.Cells.Item(132, 12).Value = 5638.9413  # L132: 5925.75 -> 5638.9413
$ws.Cells.Item(132, 13).Value = 307.7999799999998  # M132: 281.5333000000001 -> 307.7999799999998
$ws.Cells.Item(132, 14).Value = -10698.9413  # N132: -10985.75 -> -10698.9413

# Hunk 42: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1423.4166  # H136: 1238.6333 -> 1423.4166
$ws.Cells.Item(136, 9).Value = 1307.6  # I136: 1160.5834 -> 1307.6
$ws.Cells.Item(136, 10).Value = 2002.5  # J136: 1550.8334 -> 2002.5
$ws.Cells.Item(136, 11).Value = 3922.8  # K136: 3481.7502 -> 3922.8
$ws.Cells.Item(136, 12).Value = 6007.5  # L136: 4652.5002 -> 6007.5
$ws.Cells.Item(136, 13).Value = -1372.8  # M136: -931.7501999999999 -> -1372.8
$ws.Cells.Item(136, 14).Value = -11107.5  # N136: -9752.5002 -> -11107.5
